$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '79.767.00'
$ws.Range("E2").Value = '  +4.37%  '
$ws.Range("D3").Value = '3.203.12'
$ws.Range("E3").Value = '  +5.12%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '637.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.248'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +19.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +10.79%  '
$ws.Range("D10").Value = '3.201.46'
$ws.Range("E10").Value = '  +5.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.623'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +40.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000256'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +32.36%  '
$ws.Range("E13").Value = '  +3.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.69%  '
$ws.Range("D15").Value = '3.792.41'
$ws.Range("E15").Value = '  +5.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +11.51%  '
$ws.Range("D17").Value = '80.123.09'
$ws.Range("E17").Value = '  +4.94%  '
$ws.Range("D18").Value = '3.200.27'
$ws.Range("E18").Value = '  +5.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '443.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +18.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +27.76%  '
$ws.Range("E23").Value = '  +20.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +10.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '77.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.94%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +9.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.08%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("E31").Value = '  +5.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '534.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.147'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +27.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.60%  '
$ws.Range("E36").Value = '  +18.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  +6.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '192.55'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +8.46%  '
$ws.Range("E44").Value = '  +11.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.807'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.51%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.61%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +15.50%  '
$ws.Range("E50").Value = '  +5.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.65%  '
